$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
# Row 40 (G=5505)
$ws.Range("H40").Value = 5195.6665
$ws.Range("I40").Value = 4970
$ws.Range("J40").Value = 5356.857
$ws.Range("K40").Value = 4970
$ws.Range("L40").Value = 5356.857
$ws.Range("M40").Value = -4795
$ws.Range("N40").Value = -5706.857

# Row 62 (G=27781)
$ws.Range("H62").Value = 4533
$ws.Range("I62").Value = 4533
$ws.Range("K62").Value = 4533
$ws.Range("M62").Value = -3909

# Row 64 (G=5506)
$ws.Range("H64").Value = 9615.053
$ws.Range("J64").Value = 10360.5625
$ws.Range("L64").Value = 10360.5625
$ws.Range("N64").Value = -10856.5625

# Row 65 (G=27781)
$ws.Range("H65").Value = 4533
$ws.Range("I65").Value = 4533
$ws.Range("K65").Value = 22665
$ws.Range("M65").Value = -19545

# Row 67 (G=5506)
$ws.Range("H67").Value = 9615.053
$ws.Range("J67").Value = 10360.5625
$ws.Range("L67").Value = 10360.5625
$ws.Range("N67").Value = -12076.5625

# Row 70 (G=12604)
$ws.Range("H70").Value = 2460.2
$ws.Range("J70").Value = 2800
$ws.Range("L70").Value = 8400
$ws.Range("N70").Value = -8940

# Row 73 (G=12604)
$ws.Range("H73").Value = 2460.2
$ws.Range("J73").Value = 2800
$ws.Range("L73").Value = 8400
$ws.Range("N73").Value = -10272

# Row 99 (G=19883)
$ws.Range("H99").Value = 770.8333
$ws.Range("I99").Value = 318.75
$ws.Range("K99").Value = 956.25
$ws.Range("M99").Value = 541.75

# Row 132 (G=44049)
$ws.Range("H132").Value = 15966.486
$ws.Range("I132").Value = 2057.647
$ws.Range("K132").Value = 6172.941
$ws.Range("M132").Value = -3642.941

# Row 137 (G=44013)
$ws.Range("H137").Value = 4831.1943
$ws.Range("I137").Value = 5045.269
$ws.Range("J137").Value = 4274.6
$ws.Range("K137").Value = 15135.807
$ws.Range("L137").Value = 12823.8
$ws.Range("M137").Value = -12585.807
$ws.Range("N137").Value = -17923.8

# Row 141 (G=44161)
$ws.Range("H141").Value = 6207.2354
$ws.Range("I141").Value = 5247.923
$ws.Range("K141").Value = 15743.769
$ws.Range("M141").Value = -10563.769

$ws = $wb.Sheets.Item("ARM")
# Row 61 (G=43999)
$ws.Range("H61").Value = 4248.7144
$ws.Range("I61").Value = 4127
$ws.Range("J61").Value = 4467.8
$ws.Range("K61").Value = 4127
$ws.Range("L61").Value = 4467.8
$ws.Range("M61").Value = -3915
$ws.Range("N61").Value = -4891.8

# Row 101 (G=18518)
$ws.Range("H101").Value = 43733
$ws.Range("J101").Value = 43733
$ws.Range("L101").Value = 43733
$ws.Range("N101").Value = -50223

# Row 105 (G=18699)
$ws.Range("H105").Value = 10370
$ws.Range("J105").Value = 10370
$ws.Range("L105").Value = 10370
$ws.Range("N105").Value = -17358

# Row 132 (G=43997)
$ws.Range("H132").Value = 2684.25
$ws.Range("I132").Value = 2461.1904
$ws.Range("K132").Value = 7383.5712
$ws.Range("M132").Value = -4853.5712

# Row 136 (G=43999)
$ws.Range("H136").Value = 4248.7144
$ws.Range("I136").Value = 4127
$ws.Range("J136").Value = 4467.8
$ws.Range("K136").Value = 12381
$ws.Range("L136").Value = 13403.4
$ws.Range("M136").Value = -9831
$ws.Range("N136").Value = -18503.4

$ws = $wb.Sheets.Item("BSM")
# Row 82 (G=11877)
$ws.Range("H82").Value = 5038.5
$ws.Range("I82").Value = 5038.5
$ws.Range("K82").Value = 5038.5
$ws.Range("M82").Value = -4655.5

# Row 85 (G=11877)
$ws.Range("H85").Value = 5038.5
$ws.Range("I85").Value = 5038.5
$ws.Range("K85").Value = 5038.5
$ws.Range("M85").Value = -3712.5

# Row 107 (G=27706)
$ws.Range("H107").Value = 4342.689
$ws.Range("I107").Value = 3797.543
$ws.Range("J107").Value = 6250.7
$ws.Range("K107").Value = 3797.543
$ws.Range("L107").Value = 6250.7
$ws.Range("M107").Value = -1877.543
$ws.Range("N107").Value = -10090.7

# Row 134 (G=43998)
$ws.Range("H134").Value = 3739.4
$ws.Range("I134").Value = 2899.8
$ws.Range("J134").Value = 4579
$ws.Range("K134").Value = 8699.400000000001
$ws.Range("L134").Value = 13737
$ws.Range("M134").Value = -6164.400000000001
$ws.Range("N134").Value = -18807

$ws = $wb.Sheets.Item("CRP")
# Row 62 (G=12580)
$ws.Range("H62").Value = 38465244
$ws.Range("J62").Value = 55559216
$ws.Range("L62").Value = 55559216
$ws.Range("N62").Value = -55560464

# Row 65 (G=12580)
$ws.Range("H65").Value = 38465244
$ws.Range("J65").Value = 55559216
$ws.Range("L65").Value = 277796080
$ws.Range("N65").Value = -277802320

# Row 99 (G=36198)
$ws.Range("H99").Value = 11812528
$ws.Range("I99").Value = 1744645.1
$ws.Range("K99").Value = 1744645.1
$ws.Range("M99").Value = -1743147.1

# Row 122 (G=36196)
$ws.Range("H122").Value = 343326.78
$ws.Range("J122").Value = 3663
$ws.Range("L122").Value = 10989
$ws.Range("N122").Value = -15889

# Row 126 (G=36198)
$ws.Range("H126").Value = 11812528
$ws.Range("I126").Value = 1744645.1
$ws.Range("K126").Value = 5233935.300000001
$ws.Range("M126").Value = -5231465.300000001

# Row 134 (G=44020)
$ws.Range("H134").Value = 1740
$ws.Range("I134").Value = 1672.0476
$ws.Range("J134").Value = 2025.4
$ws.Range("K134").Value = 5016.142800000001
$ws.Range("L134").Value = 6076.200000000001
$ws.Range("M134").Value = -2481.142800000001
$ws.Range("N134").Value = -11146.2

$ws = $wb.Sheets.Item("CUL")
# Row 4 (G=4650)
$ws.Range("H4").Value = 77770050
$ws.Range("I4").Value = 84250470
$ws.Range("K4").Value = 252751410
$ws.Range("M4").Value = -252751298

# Row 60 (G=4750)
$ws.Range("H60").Value = 780.7059
$ws.Range("I60").Value = 186.91667
$ws.Range("K60").Value = 560.75001
$ws.Range("M60").Value = -309.75001

# Row 92 (G=19841)
$ws.Range("H92").Value = 544.125
$ws.Range("J92").Value = 543.5
$ws.Range("L92").Value = 1630.5
$ws.Range("N92").Value = -4126.5

# Row 94 (G=19811)
$ws.Range("H94").Value = 4131.6665
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352

# Row 95 (G=19838)
$ws.Range("H95").Value = 4950
$ws.Range("J95").Value = 4950
$ws.Range("L95").Value = 14850
$ws.Range("N95").Value = -18968

# Row 116 (G=27866)
$ws.Range("H116").Value = 2760
$ws.Range("J116").Value = 3325
$ws.Range("L116").Value = 9975
$ws.Range("N116").Value = -16859

# Row 132 (G=43972)
$ws.Range("H132").Value = 1815
$ws.Range("I132").Value = 1222.5
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 11002.5
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -8472.5
$ws.Range("N132").Value = -32060

$ws = $wb.Sheets.Item("GSM")
# Row 21 (G=4430)
$ws.Range("H21").Value = 257500
$ws.Range("I21").Value = 257500
$ws.Range("K21").Value = 257500
$ws.Range("M21").Value = -257327

# Row 30 (G=4430)
$ws.Range("H30").Value = 257500
$ws.Range("I30").Value = 257500
$ws.Range("K30").Value = 257500
$ws.Range("M30").Value = -257395

# Row 39 (G=18264)
$ws.Range("H39").Value = 76753.336
$ws.Range("J39").Value = 76753.336
$ws.Range("L39").Value = 76753.336
$ws.Range("N39").Value = -77817.336

# Row 70 (G=14146)
$ws.Range("H70").Value = 7022.048
$ws.Range("I70").Value = 5790
$ws.Range("J70").Value = 9024.125
$ws.Range("K70").Value = 5790
$ws.Range("L70").Value = 9024.125
$ws.Range("M70").Value = -5520
$ws.Range("N70").Value = -9564.125

# Row 73 (G=14146)
$ws.Range("H73").Value = 7022.048
$ws.Range("I73").Value = 5790
$ws.Range("J73").Value = 9024.125
$ws.Range("K73").Value = 5790
$ws.Range("L73").Value = 9024.125
$ws.Range("M73").Value = -4854
$ws.Range("N73").Value = -10896.125

# Row 97 (G=19940)
$ws.Range("H97").Value = 10481.3
$ws.Range("I97").Value = 575.4286
$ws.Range("J97").Value = 33595
$ws.Range("K97").Value = 575.4286
$ws.Range("L97").Value = 33595
$ws.Range("M97").Value = -79.42859999999996
$ws.Range("N97").Value = -34587

# Row 107 (G=27802)
$ws.Range("H107").Value = 479.35294
$ws.Range("J107").Value = 728.4
$ws.Range("L107").Value = 728.4
$ws.Range("N107").Value = -4568.4

# Row 126 (G=36184)
$ws.Range("H126").Value = 2494.3333
$ws.Range("I126").Value = 2494.3333
$ws.Range("K126").Value = 7482.999899999999
$ws.Range("M126").Value = -5012.999899999999

$ws = $wb.Sheets.Item("LTW")
# Row 46 (G=5282)
$ws.Range("H46").Value = 3044.2727
$ws.Range("I46").Value = 1366.3334
$ws.Range("K46").Value = 1366.3334
$ws.Range("M46").Value = -1178.3334

# Row 101 (G=18549)
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 132 (G=44058)
$ws.Range("H132").Value = 1420.3334
$ws.Range("I132").Value = 1015.45
$ws.Range("J132").Value = 3444.75
$ws.Range("K132").Value = 3046.35
$ws.Range("L132").Value = 10334.25
$ws.Range("M132").Value = -516.3500000000004
$ws.Range("N132").Value = -15394.25

# Row 136 (G=44060)
$ws.Range("H136").Value = 2799.8462
$ws.Range("I136").Value = 2185.1538
$ws.Range("K136").Value = 6555.4614
$ws.Range("M136").Value = -4005.4614

$ws = $wb.Sheets.Item("WVR")
# Row 98 (G=18374)
$ws.Range("H98").Value = 53331.668
$ws.Range("J98").Value = 53331.668
$ws.Range("L98").Value = 53331.668
$ws.Range("N98").Value = -59321.668

# Row 103 (G=18548)
$ws.Range("H103").Value = 45512.625
$ws.Range("J103").Value = 45512.625
$ws.Range("L103").Value = 45512.625
$ws.Range("N103").Value = -47856.625

# Row 136 (G=44031)
$ws.Range("H136").Value = 1899.2142
$ws.Range("I136").Value = 1882.4166
$ws.Range("K136").Value = 5647.2498
$ws.Range("M136").Value = -3097.2498
